$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '67.908.85'
Set-TextValue 'E2' '  +0.06%  '
Set-TextValue 'D3' '3.337.82'
Set-TextValue 'E3' '  +0.83%  '
Set-TextValue 'D4' '1.00'
Set-TextValue 'E4' '  +0.14%  '
Set-TextValue 'D5' '583.22'
Set-TextValue 'E5' '  +0.33%  '
Set-TextValue 'D6' '175.03'
Set-TextValue 'E6' '  +0.00%  '
Set-TextValue 'E7' '  +0.07%  '
Set-TextValue 'E8' '  +1.64%  '
Set-TextValue 'E9' '  +4.28%  '
Set-TextValue 'D10' '0.580'
Set-TextValue 'E10' '  +0.93%  '
Set-TextValue 'D11' '47.20'
Set-TextValue 'E11' '  +3.91%  '
Set-TextValue 'E12' '  +1.88%  '
Set-TextValue 'D13' '697.23'
Set-TextValue 'E13' '  +3.65%  '
Set-TextValue 'D14' '3.873.26'
Set-TextValue 'E14' '  +0.65%  '
Set-TextValue 'D15' '8.37'
Set-TextValue 'E15' '  +0.34%  '
Set-TextValue 'D16' '68.007.35'
Set-TextValue 'E16' '  +0.23%  '
Set-TextValue 'E17' '  +0.73%  '
Set-TextValue 'D18' '3.327.35'
Set-TextValue 'E18' '  +0.24%  '
Set-TextValue 'D19' '17.42'
Set-TextValue 'E19' '  +0.03%  '
Set-TextValue 'D20' '11.13'
Set-TextValue 'E20' '  +2.46%  '
Set-TextValue 'D21' '0.894'
Set-TextValue 'E21' '  +0.80%  '
Set-TextValue 'D22' '5.41'
Set-TextValue 'E22' '  +0.02%  '
Set-TextValue 'D23' '16.97'
Set-TextValue 'E23' '  -1.04%  '
Set-TextValue 'D24' '101.15'
Set-TextValue 'E24' '  +3.43%  '
Set-TextValue 'E25' '  +1.32%  '
Set-TextValue 'E26' '  +0.69%  '
Set-TextValue 'D27' '9.44'
Set-TextValue 'E27' '  +2.91%  '
Set-TextValue 'D28' '33.02'
Set-TextValue 'E28' '  -0.14%  '
Set-TextValue 'D29' '8.53'
Set-TextValue 'E29' '  +1.69%  '
Set-TextValue 'D30' '6.96'
Set-TextValue 'E30' '  -1.20%  '
Set-TextValue 'D31' '574.86'
Set-TextValue 'E31' '  -3.21%  '
Set-TextValue 'E32' '  +0.71%  '
Set-TextValue 'E33' '  +1.90%  '
Set-TextValue 'D34' '3.759.01'
Set-TextValue 'E34' '  +0.43%  '
Set-TextValue 'E35' '  +0.17%  '
Set-TextValue 'D36' '56.63'
Set-TextValue 'E36' '  +2.32%  '
Set-TextValue 'D37' '3.32'
Set-TextValue 'E37' '  -1.77%  '
Set-TextValue 'D38' '35.53'
Set-TextValue 'E38' '  +10.28%  '
Set-TextValue 'E39' '  +2.95%  '
Set-TextValue 'D40' '3.14'
Set-TextValue 'E40' '  +1.99%  '
Set-TextValue 'D41' '2.61'
Set-TextValue 'E41' '  -0.44%  '
Set-TextValue 'E42' '  +1.94%  '
Set-TextValue 'B43' 'ApeXProtocol'
Set-TextValue 'C43' 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue 'D43' '3.33'
Set-TextValue 'E43' '  +1.77%  '
Set-TextValue 'B44' 'TheGraph'
Set-TextValue 'C44' 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue 'D44' '0.334'
Set-TextValue 'E44' '  +1.31%  '
Set-TextValue 'E45' '  +0.30%  '
Set-TextValue 'D46' '2.64'
Set-TextValue 'E46' '  +1.70%  '
Set-TextValue 'E47' '  +1.41%  '
Set-TextValue 'E48' '  -0.43%  '
Set-TextValue 'E49' '  -1.24%  '
Set-TextValue 'D50' '130.28'
Set-TextValue 'E50' '  +1.15%  '
Set-TextValue 'D51' '2.69'
Set-TextValue 'E51' '  +4.00%  '
